# Fixes #111 - Implemented the file use model.
#
# The "file 1 use" / "file 2 use" columns (F:H) used to carry a single
# free-text "image-source" label per file row. The file-use model instead
# records, per row, the use category for each of the listed files
# (e.g. OriginalFile / PreservationMasterFile) alongside the file name -
# so the Component row's placeholder label is cleared out and the
# Sub-component row gets the real file-use triple.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 3 (Component, file_1.jpg): no file-use data yet - clear the old
# "image-source" placeholder but keep the same cell formatting across
# F3:H3 (matches the already-styled blanks on row 2, e.g. E2/F2).
$ws.Range("E2").Copy()
$ws.Range("F3:H3").PasteSpecial($xlPasteFormats)
$ws.Range("F3:H3").ClearContents()

# Row 4 (Sub-component, file_2.jpg): file-use triple - use category,
# file name, and preservation-level use category.
$ws.Range("F4").Value = "OriginalFile"
$ws.Range("G4").Value = "image.jpg"
$ws.Range("H4").Value = "PreservationMasterFile"
$ws.Range("E2").Copy()
$ws.Range("F4:H4").PasteSpecial($xlPasteFormats)

# Leave the selection on the new file-use block.
$ws.Range("F3:I3").Select()
